$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.130.00"
$ws.Range("E2").Value = "  -0.38%  "

# Row 3
$ws.Range("D3").Value = "1.657.51"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4
$ws.Range("E4").Value = "  -0.34%  "

# Row 5
$ws.Range("D5").Value = "218.13"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
$ws.Range("D6").Value = "0.5286"
$ws.Range("E6").Value = "  +1.09%  "

# Row 7
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("E8").Value = "  -2.16%  "

# Row 9
$ws.Range("D9").Value = "0.06352"
$ws.Range("E9").Value = "  +0.89%  "

# Row 10
$ws.Range("D10").Value = "20.45"
$ws.Range("E10").Value = "  -1.81%  "

# Row 11
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +0.89%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.689.80"
$ws.Range("E12").Value = "  +1.62%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.512"
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("D14").Value = "0.5494"
$ws.Range("E14").Value = "  +0.71%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8214"
$ws.Range("E15").Value = "  +0.31%  "

# Row 16
$ws.Range("D16").Value = "65.48"
$ws.Range("E16").Value = "  +1.40%  "

# Row 17
$ws.Range("D17").Value = "26.157.97"
$ws.Range("E17").Value = "  -0.35%  "

# Row 18
$ws.Range("E18").Value = "  -0.38%  "

# Row 19
$ws.Range("D19").Value = "4.578"
$ws.Range("E19").Value = "  -1.60%  "

# Row 20
$ws.Range("D20").Value = "192.61"
$ws.Range("E20").Value = "  -0.68%  "

# Row 21
$ws.Range("E21").Value = "  +0.27%  "

# Row 22
$ws.Range("D22").Value = "6.038"
$ws.Range("E22").Value = "  -0.06%  "

# Row 23
$ws.Range("E23").Value = "  -0.37%  "

# Row 24
$ws.Range("D24").Value = "141.88"
$ws.Range("E24").Value = "  +1.38%  "

# Row 25
$ws.Range("D25").Value = "0.1249"
$ws.Range("E25").Value = "  +1.40%  "

# Row 26
$ws.Range("D26").Value = "7.283"
$ws.Range("E26").Value = "  +1.71%  "

# Row 27
$ws.Range("D27").Value = "16.19"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("D28").Value = "1.443"
$ws.Range("E28").Value = "  +2.02%  "

# Row 29
$ws.Range("D29").Value = "0.05938"
$ws.Range("E29").Value = "  -3.12%  "

# Row 30
$ws.Range("D30").Value = "1.278"
$ws.Range("E30").Value = "  -0.01%  "

# Row 31
$ws.Range("D31").Value = "3.528"
$ws.Range("E31").Value = "  -1.11%  "

# Row 32
$ws.Range("D32").Value = "3.267"
$ws.Range("E32").Value = "  +0.09%  "

# Row 33
$ws.Range("D33").Value = "1.588"
$ws.Range("E33").Value = "  -2.17%  "

# Row 34
$ws.Range("D34").Value = "0.9566"
$ws.Range("E34").Value = "  -1.19%  "

# Row 35
$ws.Range("D35").Value = "2.792"
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("D36").Value = "2.412"
$ws.Range("E36").Value = "  -0.48%  "

# Row 37
$ws.Range("D37").Value = "0.5718"
$ws.Range("E37").Value = "  +0.63%  "

# Row 38
$ws.Range("D38").Value = "0.01619"
$ws.Range("E38").Value = "  +1.34%  "

# Row 39
$ws.Range("E39").Value = "  -2.81%  "

# Row 40
$ws.Range("D40").Value = "0.8470"
$ws.Range("E40").Value = "  -0.96%  "

# Row 41
$ws.Range("E41").Value = "  -0.25%  "

# Row 42
$ws.Range("D42").Value = "103.18"
$ws.Range("E42").Value = "  +2.84%  "

# Row 43
$ws.Range("D43").Value = "1.025.72"
$ws.Range("E43").Value = "  +1.07%  "

# Row 44
$ws.Range("D44").Value = "1.802.26"
$ws.Range("E44").Value = "  -0.13%  "

# Row 45
$ws.Range("D45").Value = "57.41"
$ws.Range("E45").Value = "  +0.57%  "

# Row 46
$ws.Range("E46").Value = "  -0.25%  "

# Row 47
$ws.Range("D47").Value = "1.489"
$ws.Range("E47").Value = "  +0.36%  "

# Row 48
$ws.Range("E48").Value = "  +1.77%  "

# Row 49
$ws.Range("D49").Value = "0.05157"
$ws.Range("E49").Value = "  -0.54%  "

# Row 50
$ws.Range("D50").Value = "7.815"
$ws.Range("E50").Value = "  -1.70%  "

# Row 51
$ws.Range("D51").Value = "0.09716"
$ws.Range("E51").Value = "  -0.07%  "
